$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("alpha")
$ws.Range("C2").Value = 0.004901990615312879

$ws = $wb.Worksheets.Item("beta")
$ws.Range("B2").Value = 0.000000000001888716660614536
$ws.Range("C2").Value = 0.000000001275288575208966
$ws.Range("D2").Value = 0.000000000000009114301403721972
$ws.Range("E2").Value = 0.00000000007463033791243776
$ws.Range("B3").Value = 0.0000003311150960402586
$ws.Range("C3").Value = 0.0001520139119407046
$ws.Range("D3").Value = 0.000000001597848341980173
$ws.Range("E3").Value = 0.0000204779531999528
$ws.Range("B4").Value = 0.000000000004942307295675126
$ws.Range("C4").Value = 0.000000002933562112833032
$ws.Range("D4").Value = 0.00000000000002384988667804759
$ws.Range("E4").Value = 0.0000000002312660888222127
$ws.Range("B5").Value = 0.000000000001083556910082177
$ws.Range("C5").Value = 0.0000000003921632668376349
$ws.Range("D5").Value = 0.000000000000005228875496529637
$ws.Range("E5").Value = 0.00000000007076279347027518
$ws.Range("B6").Value = 0.00000008219302300846927
$ws.Range("C6").Value = 0.00002754043675622849
$ws.Range("D6").Value = 0.0000000003966354512584736
$ws.Range("E6").Value = 0.000002566286801891313
$ws.Range("B7").Value = 0.0000000000004277407425059453
$ws.Range("C7").Value = 0.0000000004960814889487753
$ws.Range("D7").Value = 0.000000000000002064130703745967
$ws.Range("E7").Value = 0.00000000004356826553813024
$ws.Range("B8").Value = 0.00000007208104801007146
$ws.Range("C8").Value = 0.00003479992289019614
$ws.Range("D8").Value = 0.0000000003478385142460627
$ws.Range("E8").Value = 0.000004149066605309422
$ws.Range("B9").Value = 0.00005552126093801832
$ws.Range("C9").Value = 0.01580741893983122
$ws.Range("D9").Value = 0.000000267926638789295
$ws.Range("E9").Value = 0.001181082949050137
$ws.Range("B10").Value = 0.0000000009919205748422449
$ws.Range("C10").Value = 0.0000009704357130698753
$ws.Range("D10").Value = 0.000000000004786669846351541
$ws.Range("E10").Value = 0.0000001240778822929603
$ws.Range("B11").Value = 0.00000002607275002210738
$ws.Range("C11").Value = 0.00000646575192691772
$ws.Range("D11").Value = 0.0000000001258181849510797
$ws.Range("E11").Value = 0.0000007790060317776939
$ws.Range("B12").Value = 0.0000002882635963399341
$ws.Range("C12").Value = 0.00007875479822493713
$ws.Range("D12").Value = 0.000000001391061642834319
$ws.Range("E12").Value = 0.000007629205271513404
$ws.Range("B13").Value = 0.000002070227589197718
$ws.Range("C13").Value = 0.002090338270323153
$ws.Range("D13").Value = 0.000000009990211139509604
$ws.Range("E13").Value = 0.0003003534730637168
$ws.Range("B14").Value = 0.0001171965117768293
$ws.Range("C14").Value = 0.07901761536342963
$ws.Range("D14").Value = 0.0000005655503305886672
$ws.Range("E14").Value = 0.008704493875427873
$ws.Range("B15").Value = 0.00002834367799707585
$ws.Range("C15").Value = 0.03254705543766866
$ws.Range("D15").Value = 0.0000001367769075914953
$ws.Range("E15").Value = 0.003708566133128645
$ws.Range("B16").Value = 0.00000000004285651945574165
$ws.Range("C16").Value = 0.0000000541247229320031
$ws.Range("D16").Value = 0.000000000000206810922770709
$ws.Range("E16").Value = 0.000000008567366045548235

$ws = $wb.Worksheets.Item("chi")
$ws.Range("B2").Value = 0.0005464301518332064
$ws.Range("C2").Value = 0.01682290957107692
$ws.Range("D2").Value = 0.1209770084279233
$ws.Range("B3").Value = 0.001470880226390878
$ws.Range("C3").Value = 0.05253496737070348
$ws.Range("D3").Value = 0.3854214682626809
$ws.Range("B4").Value = 0.0001755330892994764
$ws.Range("C4").Value = 0.01000683835677189
$ws.Range("D4").Value = 0.1292248594463174
$ws.Range("B5").Value = 15897.37555411554
$ws.Range("C5").Value = 268494.5698219873
$ws.Range("D5").Value = 2045267.642371122
$ws.Range("B6").Value = 0.3001855257855535
$ws.Range("C6").Value = 28.20731543155558
$ws.Range("D6").Value = 121.3263390867791
$ws.Range("B7").Value = 649.4434645984869
$ws.Range("C7").Value = 69731.18964552312
$ws.Range("D7").Value = 295801.1293346361
$ws.Range("B8").Value = 9910.125433210746
$ws.Range("C8").Value = 848136.6386005477
$ws.Range("D8").Value = 5851581.655258
$ws.Range("B9").Value = 0.01307908998356739
$ws.Range("C9").Value = 1.953189490552106
$ws.Range("D9").Value = 11.13531392637941

$ws = $wb.Worksheets.Item("delta")
$ws.Range("B2").Value = 1815.876926867525
$ws.Range("C2").Value = 144759.9745489098
$ws.Range("D2").Value = 11180980.05596905
$ws.Range("B3").Value = 113.2379030376425
$ws.Range("C3").Value = 4708.095903387306
$ws.Range("D3").Value = 428311.0416143669
$ws.Range("B4").Value = 0.00041476340803681
$ws.Range("C4").Value = 0.2513557211731482
$ws.Range("D4").Value = 1.381974847338104
$ws.Range("B5").Value = 32.4006185797438
$ws.Range("C5").Value = 591.8579531633486
$ws.Range("D5").Value = 104437.3026867526
$ws.Range("B6").Value = 25.61344014379068
$ws.Range("C6").Value = 956.2294283789743
$ws.Range("D6").Value = 93348.03064671246
$ws.Range("B7").Value = 11.29892808339415
$ws.Range("C7").Value = 178.5785572526005
$ws.Range("D7").Value = 32913.35974735782
$ws.Range("B8").Value = 123.8361196996238
$ws.Range("C8").Value = 1752.850663930286
$ws.Range("D8").Value = 364768.4084493495
$ws.Range("B9").Value = 25432.21275832142
$ws.Range("C9").Value = 2005103.666620673
$ws.Range("D9").Value = 154176927.461264
